# Updated cryptos list on Mon May 27 21:26:04 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# existing coin rows, and swaps the Monero / Arweave rows (50 <-> 51) with
# their refreshed figures as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    # Some refreshed prices are plain "123.45"-style numbers which Excel's
    # auto-detection would silently coerce to a Number cell. Force text
    # storage (matching the source data's text cells) and then drop the
    # temporary text number-format so no stray style sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# --- Row 2 : Bitcoin ---
Set-TextValue 2 4 "69.689.96"
$ws.Cells.Item(2, 5).Value = "  +1.52%  "

# --- Row 3 : Ethereum ---
Set-TextValue 3 4 "3.892.79"
$ws.Cells.Item(3, 5).Value = "  +1.06%  "

# --- Row 4 : TetherUSD ---
Set-TextValue 4 4 "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# --- Row 5 : BNB ---
Set-TextValue 5 4 "605.19"
$ws.Cells.Item(5, 5).Value = "  +0.84%  "

# --- Row 6 : Solana ---
Set-TextValue 6 4 "170.80"
$ws.Cells.Item(6, 5).Value = "  +5.43%  "

# --- Row 7 : LidoStakedEther ---
Set-TextValue 7 4 "3.888.45"
$ws.Cells.Item(7, 5).Value = "  +1.05%  "

# --- Row 8 : USDC ---
$ws.Cells.Item(8, 5).Value = "  -0.30%  "

# --- Row 9 : XRP ---
$ws.Cells.Item(9, 5).Value = "  +1.07%  "

# --- Row 10 : Dogecoin ---
Set-TextValue 10 4 "0.168"
$ws.Cells.Item(10, 5).Value = "  +1.01%  "

# --- Row 11 : Toncoin ---
$ws.Cells.Item(11, 5).Value = "  +0.91%  "

# --- Row 12 : Cardano ---
$ws.Cells.Item(12, 5).Value = "  +2.19%  "

# --- Row 13 : ShibaInu ---
$ws.Cells.Item(13, 5).Value = "  +5.54%  "

# --- Row 14 : Avalanche ---
Set-TextValue 14 4 "38.25"
$ws.Cells.Item(14, 5).Value = "  +4.33%  "

# --- Row 15 : WrappedliquidstakedEther2.0 ---
Set-TextValue 15 4 "4.549.65"
$ws.Cells.Item(15, 5).Value = "  +1.16%  "

# --- Row 16 : WrappedEther ---
Set-TextValue 16 4 "3.904.27"
$ws.Cells.Item(16, 5).Value = "  +1.31%  "

# --- Row 17 : WrappedBTC ---
Set-TextValue 17 4 "69.706.60"
$ws.Cells.Item(17, 5).Value = "  +1.31%  "

# --- Row 18 : Chainlink ---
Set-TextValue 18 4 "18.79"
$ws.Cells.Item(18, 5).Value = "  +10.10%  "

# --- Row 19 : Polkadot ---
Set-TextValue 19 4 "7.64"
$ws.Cells.Item(19, 5).Value = "  +1.56%  "

# --- Row 20 : TRON ---
$ws.Cells.Item(20, 5).Value = "  -0.64%  "

# --- Row 21 : Uniswap ---
Set-TextValue 21 4 "11.16"
$ws.Cells.Item(21, 5).Value = "  -0.86%  "

# --- Row 22 : BitcoinCash ---
Set-TextValue 22 4 "490.74"
$ws.Cells.Item(22, 5).Value = "  +1.42%  "

# --- Row 23 : Polygon ---
$ws.Cells.Item(23, 5).Value = "  +4.68%  "

# --- Row 24 : PEPE ---
$ws.Cells.Item(24, 5).Value = "  +2.86%  "

# --- Row 25 : Litecoin ---
Set-TextValue 25 4 "85.31"
$ws.Cells.Item(25, 5).Value = "  +1.84%  "

# --- Row 26 ---
$ws.Cells.Item(26, 5).Value = "  +4.38%  "

# --- Row 27 ---
Set-TextValue 27 4 "12.33"
$ws.Cells.Item(27, 5).Value = "  +2.26%  "

# --- Row 28 : RenderToken ---
$ws.Cells.Item(28, 5).Value = "  +2.32%  "

# --- Row 29 : Dai ---
$ws.Cells.Item(29, 5).Value = "  +0.24%  "

# --- Row 30 : PancakeSwap ---
$ws.Cells.Item(30, 5).Value = "  +1.30%  "

# --- Row 31 : WrappedeETH ---
Set-TextValue 31 4 "4.045.05"
$ws.Cells.Item(31, 5).Value = "  +1.01%  "

# --- Row 32 : ImmutableX ---
$ws.Cells.Item(32, 5).Value = "  +2.65%  "

# --- Row 33 : NEARProtocol ---
Set-TextValue 33 4 "7.83"
$ws.Cells.Item(33, 5).Value = "  +0.10%  "

# --- Row 34 : EthereumClassic ---
Set-TextValue 34 4 "31.93"
$ws.Cells.Item(34, 5).Value = "  -0.19%  "

# --- Row 35 : RenzoRestakedETH ---
Set-TextValue 35 4 "3.861.04"
$ws.Cells.Item(35, 5).Value = "  +1.57%  "

# --- Row 36 : Hedera ---
$ws.Cells.Item(36, 5).Value = "  +0.73%  "

# --- Row 37 : Filecoin ---
Set-TextValue 37 4 "6.12"
$ws.Cells.Item(37, 5).Value = "  +4.56%  "

# --- Row 38 : Mantle ---
$ws.Cells.Item(38, 5).Value = "  +0.82%  "

# --- Row 39 : Kaspa ---
$ws.Cells.Item(39, 5).Value = "  +2.01%  "

# --- Row 40 : dogwifhat ---
Set-TextValue 40 4 "3.36"
$ws.Cells.Item(40, 5).Value = "  +13.80%  "

# --- Row 41 : FirstDigitalUSD ---
$ws.Cells.Item(41, 5).Value = "  +0.10%  "

# --- Row 42 : TheGraph ---
$ws.Cells.Item(42, 5).Value = "  +4.01%  "

# --- Row 43 : Stacks ---
$ws.Cells.Item(43, 5).Value = "  +6.52%  "

# --- Row 44 : Bittensor ---
Set-TextValue 44 4 "436.99"
$ws.Cells.Item(44, 5).Value = "  +2.18%  "

# --- Row 45 : OKB ---
Set-TextValue 45 4 "48.29"
$ws.Cells.Item(45, 5).Value = "  -0.43%  "

# --- Row 46 : Cosmos ---
Set-TextValue 46 4 "8.68"
$ws.Cells.Item(46, 5).Value = "  +3.87%  "

# --- Row 47 : USDe (unchanged) ---

# --- Row 48 : VeChain ---
$ws.Cells.Item(48, 5).Value = "  +3.34%  "

# --- Row 49 : FLOKI ---
Set-TextValue 49 4 "0.000274"
$ws.Cells.Item(49, 5).Value = "  +21.81%  "

# --- Rows 50/51 : Monero and Arweave swap ranking positions, with
#     refreshed price/volume figures ---
$ws.Cells.Item(50, 2).Value = "Arweave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue 50 4 "40.60"
$ws.Cells.Item(50, 5).Value = "  +5.47%  "

$ws.Cells.Item(51, 2).Value = "Monero"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 51 4 "143.94"
$ws.Cells.Item(51, 5).Value = "  +0.66%  "
